$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the existing row 185 so the old rows
# 184 and 185 can be relocated to rows 186 and 187 respectively, leaving
# rows 184 and 185 free for the updated / new data.
$ws.Rows.Item(186).Insert()
$ws.Rows.Item(186).Insert()

# Relocate the former row 185 (Asterix / "1a nueva(o)") down to row 187.
$ws.Rows.Item(185).Copy()
$ws.Rows.Item(187).PasteSpecial()

# Relocate the former row 184 (Asterix / "1a (cosecha lavada)") down to row 186.
$ws.Rows.Item(184).Copy()
$ws.Rows.Item(186).PasteSpecial()

$excel.CutCopyMode = 0

# Update row 184 in place with the new price data for "Región del Maule".
$ws.Cells.Item(184, 4).Value = 44628
$ws.Cells.Item(184, 10).Value = 270
$ws.Cells.Item(184, 11).Value = 9000
$ws.Cells.Item(184, 12).Value = 10000
$ws.Cells.Item(184, 13).Value = 9556
$ws.Cells.Item(184, 15).Value = "Región del Maule"
$ws.Cells.Item(184, 16).Value = 382

# Fill in the brand new row 185 for "Rosara" / "Región de Los Lagos".
$ws.Cells.Item(185, 1).Value = 11
$ws.Cells.Item(185, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(185, 3).Value = "Bíobío"
$ws.Cells.Item(185, 4).Value = 44628
$ws.Cells.Item(185, 5).Value = 8
$ws.Cells.Item(185, 6).Value = 100114001
$ws.Cells.Item(185, 7).Value = "Papa"
$ws.Cells.Item(185, 8).Value = "Rosara"
$ws.Cells.Item(185, 9).Value = "1a (cosecha)"
$ws.Cells.Item(185, 10).Value = 250
$ws.Cells.Item(185, 11).Value = 8000
$ws.Cells.Item(185, 12).Value = 8500
$ws.Cells.Item(185, 13).Value = 8200
$ws.Cells.Item(185, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(185, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(185, 16).Value = 328
$ws.Cells.Item(185, 17).Value = 25
$ws.Cells.Item(185, 18).Value = "Hortaliza"
